$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 1.649481363816475)
    3  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 9.295990156953671)
    4  = @(1.459612070389937, 0.3127903958511391, 3.900430680208489, 8.660232485948974, 14.33306563239854)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    6  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    7  = @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 3.781711156805759)
    8  = @(0.003994804209775715, 0.002777888934908601, 0.1575252929769615, 0.496779210170732, 0.6610771962923778)
    9  = @(0.3048080303191223, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.626907116734944)
    10 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 6.201049113329182)
    11 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    12 = @(0.6753301551942219, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 2.020003488691574)
    13 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
